# Insert a new weekly price record for "Bruselas (repollito)" at Vega
# Modelo de Temuco. The new observation is inserted as row 94 (existing
# rows 94..175 shift down to 95..176), matching the published source data
# being prepended in front of the most recent previously-recorded entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 94, pushing rows 94-175 down to 95-176.
$ws.Rows.Item(94).Insert()

$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 45096
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = 100112035
$ws.Range("G94").Value = "Bruselas (repollito)"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 95
$ws.Range("K94").Value = 28000
$ws.Range("L94").Value = 28000
$ws.Range("M94").Value = 28000
$ws.Range("N94").Value = "$/malla 15 kilos"
$ws.Range("O94").Value = "Región Metropolitana"
$ws.Range("P94").Value = 1867
$ws.Range("Q94").Value = 15
$ws.Range("R94").Value = "Hortaliza"
